# TASK_MP.8.xlsx edit:
#   Change the per-detector "best" summary formula from MIN(...) to MAX(...)
#   on each sheet (commit message: "Changed the min function to max function
#   for maximum total time acquired for processing"). The grand-summary cell
#   at the bottom of each sheet (which picks the best among the five/six
#   per-detector results) keeps using MIN so it still reports the best
#   (lowest) of the newly-computed per-detector maxima.
#
#   Also reproduces the cell-selection (active cell) that was left behind in
#   each sheet's view after the edits were made, and leaves "SHI TOMASI" as
#   the final active sheet/tab.

$wb = $excel.ActiveWorkbook

# Per-sheet list of cells whose formula flips from MIN(range) to MAX(range),
# and the cell that ends up selected on that sheet afterwards.
$sheetEdits = @(
    @{ Name = "SHI TOMASI"; Cells = @("F12","F24","F36","F48","F72"); Selection = "H3" },
    @{ Name = "HARRIS";     Cells = @("F12","F24","F36","F48","F72"); Selection = "F72" },
    @{ Name = "FAST";       Cells = @("F12","F24","F36","F48","F72"); Selection = "F76" },
    @{ Name = "BRISK";      Cells = @("F12","F24","F36","F48","F72"); Selection = "F36" },
    @{ Name = "ORB";        Cells = @("F12","F24","F36","F48","F72"); Selection = "F12" },
    @{ Name = "AKAZE";      Cells = @("F12","F24","F36","F48","F60","F73"); Selection = "F12" },
    @{ Name = "SIFT";       Cells = @("F12","F24","F48","F73"); Selection = "F73" }
)

foreach ($edit in $sheetEdits) {
    $ws = $wb.Worksheets.Item($edit.Name)
    $ws.Activate()

    foreach ($cellRef in $edit.Cells) {
        $cell = $ws.Range($cellRef)
        $formula = $cell.Formula
        $newFormula = $formula -replace "^=MIN\(", "=MAX("
        $cell.Formula = $newFormula
    }

    $ws.Range($edit.Selection).Select()
}

# Re-activate the first sheet last so it ends up the selected/visible tab.
$firstSheet = $wb.Worksheets.Item("SHI TOMASI")
$firstSheet.Activate()
$firstSheet.Range("H3").Select()
